# Replace the visible text of a single paragraph while leaving every other
# run (including the leading empty <w:r/>) exactly as-is.
#
# A plain Find.Execute / Range.Text replace normalizes (merges) the whole
# paragraph's runs, which silently drops the empty <w:r/> that precedes the
# text run in this document. Round-tripping through Range.WordOpenXML /
# Range.InsertXML instead only rewrites the <w:t> payload inside the
# exported fragment, so the sibling runs survive untouched. InsertXML's
# exporter stamps placeholder w14:paraId/w14:textId/w:rsidR/w:rsidRDefault
# attributes onto the <w:p> it hands back - strip those back off so the
# paragraph tag is saved exactly as clean as it started.
function Replace-ParagraphText {
    param($doc, $index, $oldText, $newText)

    $para = $doc.Paragraphs.Item($index)
    $range = $para.Range

    $xml = $range.WordOpenXML

    if (-not $xml.Contains($oldText)) {
        throw "Replace-ParagraphText: paragraph $index does not contain expected text: $oldText"
    }

    $xml = $xml.Replace($oldText, $newText)

    $stampPattern = ' w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+" w:rsidR="[0-9A-Fa-f]+" w:rsidRDefault="[0-9A-Fa-f]+"'
    $xml = [System.Text.RegularExpressions.Regex]::Replace($xml, $stampPattern, "")

    $range.InsertXML($xml)
}

$d = $word.ActiveDocument

# "What we like" bullet list
Replace-ParagraphText $d 35 "Beautiful and immersive graphics" "Unique gameplay mechanics with cluster-based payouts"
Replace-ParagraphText $d 36 "Thematic symbol design" "Three bonus rounds set in different time periods"
Replace-ParagraphText $d 37 "Three bonus rounds" "Thematic and beautifully designed symbol animations"
Replace-ParagraphText $d 38 "Demo mode available" "Range of betting options and demo mode for free play"

# "What we don't like" bullet list
Replace-ParagraphText $d 40 "Betting range is limited" "Limited betting range"
Replace-ParagraphText $d 41 "No progressive jackpot" "No progressive jackpot feature"

# Closing italic summary paragraph
Replace-ParagraphText $d 43 `
    "Read our review of Miles Bellhouse and the Gears of Time and play for free. Enjoy beautiful graphics, bonus rounds, and payouts up to 2000 times your bet." `
    "Read our review of Miles Bellhouse and the Gears of Time and play this slot game for free."
